# Refresh the cryptos list (Coin / Link / Price / Volume(1h)) with the
# latest scrape values (GitHub Actions run on Tue Mar 28 08:50:26 UTC 2023).
#
# Most rows only get new Price (col D) / Volume(1h) (col E) figures. Two
# rows (37 and 38) change rank and swap places: WEMIXTOKEN drops out of
# rank 35 and is replaced there by Aptos (whose price/volume also ticked),
# with WEMIXTOKEN moving down to rank 36.
#
# Price/Volume are stored as plain text in this sheet (not numbers), so any
# new value that merely *looks* numeric (e.g. "310.68") has to be forced to
# text - otherwise Excel would silently reinterpret it as a numeric cell.
function Set-TextCell($ws, $addr, $val) {
    if ($val -match '^-?\d+(\.\d+)?$') {
        # Numeric-looking string: a leading quote-prefix forces Excel to
        # store it as text instead of parsing it as a number; then drop
        # back to the default "Normal" cell style so no stray text number
        # format is left behind on the cell.
        $ws.Range($addr).Value = "'" + $val
        $ws.Range($addr).Style = "Normal"
    } else {
        $ws.Range($addr).Value = $val
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '27.018.89'),
    @('E2', '  -2.99%  '),
    @('D3', '1.727.52'),
    @('E3', '  -1.65%  '),
    @('E4', '  -0.05%  '),
    @('D5', '310.68'),
    @('D6', '1.001'),
    @('E6', '  +0.04%  '),
    @('D7', '0.4853'),
    @('E7', '  +3.52%  '),
    @('D8', '0.3488'),
    @('E8', '  -0.31%  '),
    @('D9', '43.47'),
    @('E9', '  +3.57%  '),
    @('D10', '0.07244'),
    @('E10', '  -1.55%  '),
    @('D11', '1.054'),
    @('E11', '  -2.52%  '),
    @('D12', '1.001'),
    @('E12', '  +0.00%  '),
    @('D13', '19.97'),
    @('E13', '  -2.56%  '),
    @('D14', '5.882'),
    @('E14', '  -1.65%  '),
    @('D15', '1.728.18'),
    @('E15', '  -1.53%  '),
    @('D16', '6.856'),
    @('E16', '  -4.06%  '),
    @('D17', '86.98'),
    @('E17', '  -5.43%  '),
    @('D18', '0.00001034'),
    @('E18', '  -1.84%  '),
    @('D19', '0.06394'),
    @('E19', '  -0.17%  '),
    @('E20', '  +0.05%  '),
    @('D21', '16.62'),
    @('E21', '  -0.84%  '),
    @('D22', '5.716'),
    @('E22', '  -0.56%  '),
    @('D23', '27.077.30'),
    @('E23', '  -2.86%  '),
    @('D24', '10.95'),
    @('E24', '  -1.76%  '),
    @('E25', '  -3.53%  '),
    @('D26', '154.31'),
    @('E26', '  -4.40%  '),
    @('D27', '20.01'),
    @('E27', '  +0.04%  '),
    @('D28', '1.922.69'),
    @('E28', '  -1.75%  '),
    @('D29', '2.078'),
    @('E29', '  -3.28%  '),
    @('D30', '120.95'),
    @('E30', '  -1.31%  '),
    @('D31', '1.046'),
    @('E31', '  -1.97%  '),
    @('D32', '0.09337'),
    @('E32', '  +0.04%  '),
    @('E33', '  -0.09%  '),
    @('D34', '5.400'),
    @('E34', '  -2.52%  '),
    @('D35', '0.05929'),
    @('E35', '  -2.18%  '),
    @('D36', '0.02186'),
    @('E36', '  -3.40%  '),
    @('B37', 'Aptos'),
    @('C37', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'),
    @('D37', '10.98'),
    @('E37', '  -5.43%  '),
    @('B38', 'WEMIXTOKEN'),
    @('C38', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'),
    @('D38', '1.426'),
    @('E38', '  +5.60%  '),
    @('D39', '0.1996'),
    @('E39', '  -3.25%  '),
    @('D40', '4.763'),
    @('E40', '  -2.59%  '),
    @('D41', '0.9998'),
    @('E41', '  +0.02%  '),
    @('D42', '0.5993'),
    @('E42', '  -2.23%  '),
    @('D43', '1.119'),
    @('E43', '  -4.94%  '),
    @('D44', '7.547'),
    @('E44', '  -2.76%  '),
    @('D45', '12.75'),
    @('E45', '  -2.53%  '),
    @('D46', '3.583'),
    @('E46', '  -4.00%  '),
    @('D47', '0.5629'),
    @('E47', '  -2.45%  '),
    @('D48', '118.94'),
    @('E48', '  -2.93%  '),
    @('D49', '1.849'),
    @('E49', '  -3.73%  '),
    @('D50', '1.107'),
    @('E50', '  -1.10%  '),
    @('D51', '0.06643'),
    @('E51', '  -2.23%  ')
)

foreach ($u in $updates) {
    Set-TextCell $ws $u[0] $u[1]
}
